# Apply a cyclic rotation of the species-observation data held in rows 4-6
# of the "Artfynd" sheet: the contents that used to live in row 5 now
# belong to row 4, what used to be in row 6 now belongs to row 5, and what
# used to be in row 4 now belongs to row 6 (A5->A4, A6->A5, A4->A6, etc.)
# This mirrors the reordering described in the diff for columns
# A, B, D, E, F, G, H, I, J, Q and R (the remaining columns on those rows
# happen to already be identical across rows 4-6, so no further action is
# required for them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","I","J","Q","R")

# Column I ("Antal") is stored as text (e.g. "1", "5") rather than as a
# number, so force it to keep a text format before writing back into it.
$ws.Range("I4:I6").NumberFormat = "@"

# Snapshot the current ("before") values of the three rows for every
# column that participates in the rotation.
$row4 = @{}
$row5 = @{}
$row6 = @{}

foreach ($col in $cols) {
    $row4[$col] = $ws.Range("${col}4").Value2
    $row5[$col] = $ws.Range("${col}5").Value2
    $row6[$col] = $ws.Range("${col}6").Value2
}

# Write the rotated values back: row4 receives what used to be in row5,
# row5 receives what used to be in row6, and row6 receives what used to
# be in row4.
foreach ($col in $cols) {
    $ws.Range("${col}4").Value2 = $row5[$col]
    $ws.Range("${col}5").Value2 = $row6[$col]
    $ws.Range("${col}6").Value2 = $row4[$col]
}
